$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.937.82"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.748.05"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.25"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.16"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  -20.27%  "
$ws.Range("D13").Value = "3.230.25"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.50"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "63.597.95"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "2.750.69"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.14"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.72"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.536"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.02"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.39"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "0.0₃0891"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.93"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.94"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.88"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.12"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.85"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.976"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.16"
$ws.Range("E39").Value = "  +8.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.13"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "325.58"
$ws.Range("E41").Value = "  -6.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.94"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.25"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0587"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.38"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0253"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.88"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.624"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +0.64%  "
